$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 83.57143000000001
$ws.Range("I6").Value = 83.57143000000001
$ws.Range("K6").Value = 250.71429
$ws.Range("M6").Value = -138.71429

$ws.Range("H9").Value = 171.4
$ws.Range("J9").Value = 165
$ws.Range("L9").Value = 165
$ws.Range("N9").Value = -503

$ws.Range("H76").Value = 13141.167
$ws.Range("I76").Value = 5711.75
$ws.Range("K76").Value = 5711.75
$ws.Range("M76").Value = -5396.75

$ws.Range("H79").Value = 13141.167
$ws.Range("I79").Value = 5711.75
$ws.Range("K79").Value = 5711.75
$ws.Range("M79").Value = -4619.75

$ws.Range("H92").Value = 1072
$ws.Range("I92").Value = 628.2857
$ws.Range("J92").Value = 1848.5
$ws.Range("K92").Value = 628.2857
$ws.Range("L92").Value = 1848.5
$ws.Range("M92").Value = 619.7143
$ws.Range("N92").Value = -4344.5

$ws.Range("H98").Value = 2173.3333
$ws.Range("I98").Value = 2046.2307
$ws.Range("J98").Value = 2999.5
$ws.Range("K98").Value = 2046.2307
$ws.Range("L98").Value = 2999.5
$ws.Range("M98").Value = -548.2307000000001
$ws.Range("N98").Value = -5995.5

$ws.Range("H115").Value = 4586
$ws.Range("J115").Value = 9998
$ws.Range("L115").Value = 29994
$ws.Range("N115").Value = -33128

$ws.Range("H122").Value = 2173.3333
$ws.Range("I122").Value = 2046.2307
$ws.Range("J122").Value = 2999.5
$ws.Range("K122").Value = 6138.6921
$ws.Range("L122").Value = 8998.5
$ws.Range("M122").Value = -3688.6921
$ws.Range("N122").Value = -13898.5

$ws.Range("H137").Value = 23811324
$ws.Range("I137").Value = 26317672
$ws.Range("J137").Value = 1001.5
$ws.Range("K137").Value = 78953016
$ws.Range("L137").Value = 3004.5
$ws.Range("M137").Value = -78950466
$ws.Range("N137").Value = -8104.5

$ws.Range("H138").Value = 2303.3333
$ws.Range("I138").Value = 2179.4285
$ws.Range("J138").Value = 2328.1143
$ws.Range("K138").Value = 6538.2855
$ws.Range("L138").Value = 6984.342900000001
$ws.Range("M138").Value = -1398.2855
$ws.Range("N138").Value = -17264.3429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4226.25
$ws.Range("I45").Value = 4398.6665
$ws.Range("K45").Value = 4398.6665
$ws.Range("M45").Value = -4021.6665

$ws.Range("H61").Value = 4502.5
$ws.Range("I61").Value = 2670
$ws.Range("J61").Value = 10000
$ws.Range("K61").Value = 2670
$ws.Range("L61").Value = 10000
$ws.Range("M61").Value = -2458
$ws.Range("N61").Value = -10424

$ws.Range("H75").Value = 150000
$ws.Range("J75").Value = 150000
$ws.Range("L75").Value = 150000
$ws.Range("N75").Value = -151748

$ws.Range("H78").Value = 150000
$ws.Range("J78").Value = 150000
$ws.Range("L78").Value = 450000
$ws.Range("N78").Value = -458736

$ws.Range("H122").Value = 3853.6086
$ws.Range("I122").Value = 3628.3157
$ws.Range("K122").Value = 10884.9471
$ws.Range("M122").Value = -8434.947100000001

$ws.Range("H132").Value = 7300.56
$ws.Range("I132").Value = 6118.3184
$ws.Range("J132").Value = 15970.333
$ws.Range("K132").Value = 18354.9552
$ws.Range("L132").Value = 47910.999
$ws.Range("M132").Value = -15824.9552
$ws.Range("N132").Value = -52970.999

$ws.Range("H136").Value = 4502.5
$ws.Range("I136").Value = 2670
$ws.Range("J136").Value = 10000
$ws.Range("K136").Value = 8010
$ws.Range("L136").Value = 30000
$ws.Range("M136").Value = -5460
$ws.Range("N136").Value = -35100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()

$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

$ws.Range("H134").Value = 5603.1113
$ws.Range("I134").Value = 4346.857
$ws.Range("J134").Value = 10000
$ws.Range("K134").Value = 13040.571
$ws.Range("L134").Value = 30000
$ws.Range("M134").Value = -10505.571
$ws.Range("N134").Value = -35070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3790.0667
$ws.Range("I58").Value = 2038.6666
$ws.Range("J58").Value = 4957.6665
$ws.Range("K58").Value = 2038.6666
$ws.Range("L58").Value = 4957.6665
$ws.Range("M58").Value = -1835.6666
$ws.Range("N58").Value = -5363.6665

$ws.Range("H86").Value = 55565304
$ws.Range("I86").Value = 83341490
$ws.Range("K86").Value = 83341490
$ws.Range("M86").Value = -83340367

$ws.Range("H89").Value = 55565304
$ws.Range("I89").Value = 83341490
$ws.Range("K89").Value = 416707450
$ws.Range("M89").Value = -416701834

$ws.Range("H105").Value = 5714.524
$ws.Range("I105").Value = 8331.923000000001
$ws.Range("K105").Value = 8331.923000000001
$ws.Range("M105").Value = -6584.923000000001

$ws.Range("H112").Value = 87500
$ws.Range("J112").Value = 87500
$ws.Range("L112").Value = 87500
$ws.Range("N112").Value = -90454

$ws.Range("H134").Value = 2780.8
$ws.Range("I134").Value = 2656.3572
$ws.Range("K134").Value = 7969.071599999999
$ws.Range("M134").Value = -5434.071599999999

$ws.Range("H136").Value = 3790.0667
$ws.Range("I136").Value = 2038.6666
$ws.Range("J136").Value = 4957.6665
$ws.Range("K136").Value = 6115.9998
$ws.Range("L136").Value = 14872.9995
$ws.Range("M136").Value = -3565.9998
$ws.Range("N136").Value = -19972.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 323.14285
$ws.Range("I33").Value = 150.6
$ws.Range("J33").Value = 419
$ws.Range("K33").Value = 903.5999999999999
$ws.Range("L33").Value = 2514
$ws.Range("M33").Value = -620.5999999999999
$ws.Range("N33").Value = -3080

$ws.Range("H107").Value = 600.7778
$ws.Range("J107").Value = 643.2273
$ws.Range("L107").Value = 1929.6819
$ws.Range("N107").Value = -5769.6819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1616.5454
$ws.Range("J113").Value = 1681.5
$ws.Range("L113").Value = 1681.5
$ws.Range("N113").Value = -6021.5

$ws.Range("H132").Value = 3595.2104
$ws.Range("I132").Value = 3572.7222
$ws.Range("K132").Value = 10718.1666
$ws.Range("M132").Value = -8188.1666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 5000
$ws.Range("J45").Value = 5000
$ws.Range("L45").Value = 5000
$ws.Range("N45").Value = -5814

$ws.Range("H55").Value = 619.7
$ws.Range("I55").Value = 291.41666
$ws.Range("K55").Value = 291.41666
$ws.Range("M55").Value = -118.41666

$ws.Range("H132").Value = 2293.5
$ws.Range("I132").Value = 2180.7778
$ws.Range("K132").Value = 6542.3334
$ws.Range("M132").Value = -4012.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 37698.8
$ws.Range("I34").Value = 22123.75
$ws.Range("K34").Value = 22123.75
$ws.Range("M34").Value = -21920.75

$ws.Range("H37").Value = 45000
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()

$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

$ws.Range("H62").Value = 121664.6
$ws.Range("I62").Value = 167442.58
$ws.Range("J62").Value = 14849.333
$ws.Range("K62").Value = 167442.58
$ws.Range("L62").Value = 14849.333
$ws.Range("M62").Value = -166818.58
$ws.Range("N62").Value = -16097.333

$ws.Range("H65").Value = 121664.6
$ws.Range("I65").Value = 167442.58
$ws.Range("J65").Value = 14849.333
$ws.Range("K65").Value = 837212.8999999999
$ws.Range("L65").Value = 74246.66500000001
$ws.Range("M65").Value = -834092.8999999999
$ws.Range("N65").Value = -80486.66500000001

$ws.Range("H112").Value = 30200
$ws.Range("J112").Value = 30200
$ws.Range("L112").Value = 30200
$ws.Range("N112").Value = -33154

$ws.Range("H132").Value = 1692.4667
$ws.Range("I132").Value = 1984
$ws.Range("K132").Value = 5952
$ws.Range("M132").Value = -3422
